$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.458.92"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.906.89"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.705"
$ws.Range("E5").Value = "  +10.94%  "
$ws.Range("D6").Value = "'246.69"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'40.76"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").Value = "'52.66"
$ws.Range("E10").Value = "  +8.26%  "
$ws.Range("D11").Value = "'0.0726"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "'0.0989"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "2.182.98"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'12.58"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "'0.713"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").Value = "1.909.98"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "'4.90"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "35.446.33"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'73.12"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'241.99"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'12.88"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "'5.06"
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  +5.18%  "
$ws.Range("D27").Value = "'169.31"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").Value = "'18.92"
$ws.Range("E29").Value = "  +5.33%  "
$ws.Range("E30").Value = "  +5.37%  "
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.19"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.88"
$ws.Range("E35").Value = "  +6.63%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  +9.55%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "'96.35"
$ws.Range("E40").Value = "  +5.88%  "
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'16.56"
$ws.Range("E42").Value = "  +5.46%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.0652"
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").Value = "1.355.57"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("B47").Value = "MultiversX"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D47").Value = "'46.16"
$ws.Range("E47").Value = "  -8.75%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.41"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.79"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "'12.23"
$ws.Range("E50").Value = "  -4.80%  "
$ws.Range("D51").Value = "'6.52"
$ws.Range("E51").Value = "  -1.75%  "
